# Fixing some issues with csv cascading select example and dropdowns
#
# The "survey" sheet wraps the "name" text question in its own
# begin screen / end screen block (rows 13-16):
#   13: begin screen
#   14: text | name | Enter your name | Please use your full name
#   15: note | {{#if name}}Hello {{name}}{{/if}}
#   16: end screen
#
# That screen grouping is removed so the question just flows with the
# rest of the form (matching the other examples). Deleting row 13
# ("begin screen") shifts everything up by one, so the former row 16
# ("end screen") becomes row 15 - delete that too.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

$ws.Rows.Item(13).Delete()
$ws.Rows.Item(15).Delete()
